$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E5: 21 -> 22
$ws.Range("E5").Value = 22

# F11: 6 -> 7
$ws.Range("F11").Value = 7

# H11: 6 -> 7
$ws.Range("H11").Value = 7

# E16: 254 -> 256
$ws.Range("E16").Value = 256
